# Auto-generated edit script: reorders country rows per updated ranking
# and refreshes a handful of case counts plus the "updated at" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("China", 80967, 39, 71150, 6569, 2136, 3, 3248),
    @("Italia", 41035, 0, 4440, 33190, 2498, 0, 3405),
    @("España", 20412, 2335, 1588, 17783, 939, 210, 1041),
    @("Iran", 19644, 1237, 6745, 11466, 0, 149, 1433),
    @("Alemania", 17776, 2456, 180, 17548, 2, 4, 48),
    @("Estados Unidos", 14549, 760, 125, 14206, 64, 11, 218),
    @("Francia", 10995, 0, 1295, 9328, 1122, 0, 372),
    @("Corea del Sur", 8652, 87, 2233, 6325, 59, 3, 94),
    @("Suiza", 4905, 683, 15, 4840, 0, 7, 50),
    @("Reino Unido", 3269, 0, 65, 3059, 20, 1, 145),
    @("Paises Bajos", 2994, 534, 2, 2886, 45, 30, 106),
    @("Austria", 2388, 209, 9, 2373, 14, 0, 6),
    @("Belgica", 2257, 462, 204, 2016, 164, 16, 37),
    @("Noruega", 1881, 91, 1, 1873, 27, 0, 7),
    @("Suecia", 1639, 200, 16, 1607, 21, 5, 16),
    @("Dinamarca", 1226, 75, 1, 1216, 37, 3, 9),
    @("Malasia", 1030, 130, 87, 940, 26, 1, 3),
    @("Portugal", 1020, 234, 5, 1009, 26, 2, 6),
    @("Japon", 963, 20, 215, 715, 50, 0, 33),
    @("Australia", 876, 120, 46, 823, 2, 0, 7),
    @("Canada", 873, 0, 11, 850, 1, 0, 12),
    @("Chequia", 774, 80, 4, 770, 6, 0, 0),
    @("Crucero", 712, 0, 527, 178, 14, 0, 7),
    @("Israel", 705, 28, 15, 690, 10, 0, 0),
    @("Brasil", 651, 11, 2, 642, 18, 0, 7),
    @("Irlanda", 557, 0, 5, 549, 6, 0, 3),
    @("Luxemburgo", 484, 149, 6, 473, 1, 1, 5),
    @("Pakistan", 478, 24, 13, 462, 0, 1, 3),
    @("Grecia", 464, 0, 19, 439, 16, 0, 6),
    @("Catar", 460, 0, 10, 450, 6, 0, 0),
    @("Finlandia", 450, 50, 10, 440, 2, 0, 0),
    @("Chile", 434, 92, 6, 428, 7, 0, 0),
    @("Islandia", 409, 79, 5, 404, 1, 0, 0),
    @("Singapur", 385, 40, 131, 254, 14, 0, 0),
    @("Polonia", 378, 23, 13, 359, 3, 1, 6),
    @("Indonesia", 369, 60, 17, 320, 0, 7, 32),
    @("Turquia", 359, 0, 0, 355, 0, 0, 4),
    @("Eslovenia", 341, 22, 0, 340, 6, 0, 1),
    @("Tailandia", 322, 50, 42, 279, 1, 0, 1),
    @("Rumania", 308, 31, 31, 277, 11, 0, 0),
    @("Barein", 284, 5, 110, 173, 4, 0, 1),
    @("Estonia", 283, 16, 1, 282, 1, 0, 0),
    @("Arabia Saudita", 274, 0, 8, 266, 0, 0, 0),
    @("Ecuador", 260, 0, 1, 255, 2, 1, 4),
    @("Egipto", 256, 0, 42, 207, 0, 0, 7),
    @("Hong Kong", 256, 48, 98, 154, 4, 0, 4),
    @("Rusia", 253, 54, 12, 240, 0, 0, 1),
    @("Peru", 234, 0, 1, 230, 7, 2, 3),
    @("Filipinas", 230, 13, 8, 204, 1, 1, 18),
    @("India", 223, 29, 23, 195, 0, 1, 5),
    @("Sudafrica", 202, 52, 0, 202, 0, 0, 0),
    @("Irak", 192, 0, 49, 130, 0, 0, 13),
    @("Mexico", 164, 46, 4, 159, 1, 0, 1),
    @("Libano", 163, 6, 4, 155, 3, 0, 4),
    @("Kuwait", 159, 11, 22, 137, 5, 0, 0),
    @("San Marino", 144, 0, 4, 126, 12, 0, 14),
    @("Emiratos Arabes Unidos", 140, 0, 31, 109, 2, 0, 0),
    @("Eslovaquia", 137, 13, 0, 137, 2, 0, 0),
    @("Panama", 137, 0, 1, 135, 7, 0, 1),
    @("Armenia", 136, 14, 1, 135, 2, 0, 0),
    @("Taiwan", 135, 27, 28, 105, 0, 1, 2),
    @("Bulgaria", 129, 22, 1, 125, 0, 0, 3),
    @("Colombia", 128, 20, 1, 127, 0, 0, 0),
    @("Argentina", 128, 0, 3, 122, 0, 0, 3),
    @("Serbia", 118, 15, 2, 116, 4, 0, 0),
    @("Croacia", 113, 3, 5, 107, 0, 0, 1),
    @("Letonia", 111, 25, 1, 110, 0, 0, 0),
    @("Uruguay", 94, 15, 0, 94, 0, 0, 0),
    @("Vietnam", 91, 6, 17, 74, 0, 0, 0),
    @("Argelia", 90, 0, 32, 47, 0, 2, 11),
    @("Costa Rica", 89, 2, 0, 87, 2, 1, 2),
    @("Hungria", 85, 12, 7, 75, 6, 2, 3),
    @("Islas Feroe", 80, 8, 3, 77, 0, 0, 0),
    @("Brunei", 78, 5, 1, 77, 2, 0, 0),
    @("Principado de Andorra", 75, 1, 1, 74, 2, 0, 0),
    @("Republica Dominicana", 72, 38, 0, 70, 0, 0, 2),
    @("Sri Lanka", 71, 11, 3, 68, 0, 0, 0),
    @("Republica de Macedonia", 70, 20, 1, 69, 1, 0, 0),
    @("Albania", 70, 6, 0, 68, 2, 0, 2),
    @("Jordania", 69, 0, 1, 68, 0, 0, 0),
    @("Bosnia y Herzegovina", 69, 5, 2, 67, 1, 0, 0),
    @("Bielorrusia", 69, 18, 15, 54, 0, 0, 0),
    @("Republica de Chipre", 67, 0, 0, 67, 1, 0, 0),
    @("Marruecos", 66, 3, 2, 61, 1, 1, 3),
    @("Malta", 64, 11, 2, 62, 1, 0, 0),
    @("Tunez", 54, 15, 1, 52, 2, 0, 1),
    @("Camboya", 51, 14, 1, 50, 0, 0, 0),
    @("Kazajistan", 49, 5, 0, 49, 0, 0, 0),
    @("Lituania", 49, 1, 1, 48, 1, 0, 0),
    @("Moldavia", 49, 0, 1, 47, 3, 0, 1),
    @("Oman", 48, 0, 13, 35, 0, 0, 0),
    @("Estado de Palestina", 48, 1, 17, 31, 0, 0, 0),
    @("Guadalupe", 45, 12, 0, 45, 0, 0, 0),
    @("Georgia", 44, 4, 1, 43, 1, 0, 0),
    @("Azerbaiyan", 44, 0, 7, 36, 0, 0, 1),
    @("Venezuela", 42, 0, 0, 42, 0, 0, 0),
    @("Burkina Faso", 40, 7, 4, 35, 0, 0, 1),
    @("Nueva Zelanda", 39, 11, 0, 39, 0, 0, 0),
    @("Senegal", 38, 2, 2, 36, 0, 0, 0),
    @("Uzbekistan", 33, 10, 0, 33, 0, 0, 0),
    @("Martinica", 32, 9, 0, 31, 7, 0, 1),
    @("Reunion", 28, 0, 0, 28, 0, 0, 0),
    @("Liechtenstein", 28, 0, 0, 28, 0, 0, 0),
    @("Ucrania", 26, 0, 1, 22, 0, 0, 3),
    @("Honduras", 24, 12, 0, 24, 0, 0, 0),
    @("Afganistan", 24, 2, 1, 23, 0, 0, 0),
    @("Camerun", 20, 7, 2, 18, 0, 0, 0),
    @("Banglades", 20, 2, 3, 16, 1, 0, 1),
    @("Consejo Danes para los Refugiados", 18, 4, 0, 18, 0, 0, 0),
    @("Macao", 17, 0, 10, 7, 0, 0, 0),
    @("Bolivia", 16, 1, 0, 16, 0, 0, 0),
    @("Cuba", 16, 5, 0, 15, 0, 0, 1),
    @("Jamaica", 16, 1, 2, 13, 0, 0, 1),
    @("Ghana", 16, 5, 8, 8, 0, 0, 0),
    @("Guayana Francesa", 15, 0, 0, 15, 0, 0, 0),
    @("Guyana", 15, 10, 0, 14, 0, 0, 1),
    @("Guam", 14, 2, 0, 14, 0, 0, 0),
    @("Paraguay", 13, 0, 0, 13, 1, 0, 0),
    @("Maldivas", 13, 0, 0, 13, 0, 0, 0),
    @("Montenegro", 13, 0, 0, 13, 0, 0, 0),
    @("Mauricio", 12, 5, 0, 12, 0, 0, 0),
    @("Nigeria", 12, 0, 1, 11, 0, 0, 0),
    @("Monaco", 11, 1, 0, 11, 0, 0, 0),
    @("Ruanda", 11, 0, 0, 11, 0, 0, 0),
    @("Polinesia Francesa", 11, 5, 0, 11, 0, 0, 0),
    @("Gibraltar", 10, 0, 2, 8, 0, 0, 0),
    @("Togo", 9, 8, 0, 9, 0, 0, 0),
    @("Etiopia", 9, 2, 0, 9, 0, 0, 0),
    @("Trinidad yTobago", 9, 0, 0, 9, 0, 0, 0),
    @("Costa de Marfil", 9, 0, 1, 8, 0, 0, 0),
    @("Guatemala", 9, 0, 0, 8, 0, 0, 1),
    @("Puerto Rico", 8, 2, 0, 8, 0, 0, 0),
    @("Seychelles", 7, 1, 0, 7, 0, 0, 0),
    @("Kenia", 7, 0, 0, 7, 0, 0, 0),
    @("Guinea Ecuatorial", 6, 0, 0, 6, 0, 0, 0),
    @("Mongolia", 6, 0, 0, 6, 0, 0, 0),
    @("Mayotte", 6, 2, 0, 6, 0, 0, 0),
    @("Kirguistan", 6, 3, 0, 6, 0, 0, 0),
    @("Tanzania", 6, 0, 0, 6, 0, 0, 0),
    @("Barbados", 5, 0, 0, 5, 0, 0, 0),
    @("Aruba", 5, 0, 1, 4, 0, 0, 0),
    @("San Martin (Parte Francesa)", 4, 1, 0, 4, 0, 0, 0),
    @("Surinam", 4, 3, 0, 4, 0, 0, 0),
    @("Islas Virgenes de los Estados Unidos", 3, 0, 0, 3, 0, 0, 0),
    @("San Bartolome", 3, 0, 0, 3, 0, 0, 0),
    @("Congo", 3, 0, 0, 3, 0, 0, 0),
    @("Bahamas", 3, 0, 0, 3, 0, 0, 0),
    @("Namibia", 3, 0, 0, 3, 0, 0, 0),
    @("Republica de Africa Central", 3, 2, 0, 3, 0, 0, 0),
    @("Curazao", 3, 0, 0, 2, 0, 0, 1),
    @("Gabon", 3, 0, 0, 2, 0, 1, 1),
    @("Islas Caimanes", 3, 0, 0, 2, 0, 0, 1),
    @("Butan", 2, 1, 0, 2, 0, 0, 0),
    @("Liberia", 2, 0, 0, 2, 0, 0, 0),
    @("Groenlandia", 2, 0, 0, 2, 0, 0, 0),
    @("Santa Lucia", 2, 0, 0, 2, 0, 0, 0),
    @("Guinea", 2, 1, 0, 2, 0, 0, 0),
    @("Nueva Caledonia", 2, 0, 0, 2, 0, 0, 0),
    @("Bermudas", 2, 0, 0, 2, 0, 0, 0),
    @("Benin", 2, 0, 0, 2, 0, 0, 0),
    @("Haiti", 2, 2, 0, 2, 0, 0, 0),
    @("Mauritania", 2, 0, 0, 2, 0, 0, 0),
    @("Zambia", 2, 0, 0, 2, 0, 0, 0),
    @("Sudan", 2, 0, 0, 1, 0, 0, 1),
    @("Angola", 1, 1, 0, 1, 0, 0, 0),
    @("Niger", 1, 0, 0, 1, 0, 0, 0),
    @("Isla de Man", 1, 0, 0, 1, 0, 0, 0),
    @("San Martin (Parte Holandesa)", 1, 0, 0, 1, 0, 0, 0),
    @("Antigua y Barbuda", 1, 0, 0, 1, 0, 0, 0),
    @("El Salvador", 1, 0, 0, 1, 0, 0, 0),
    @("Nicaragua", 1, 0, 0, 1, 0, 0, 0),
    @("Republica del Chad", 1, 0, 0, 1, 0, 0, 0),
    @("Suazilandia", 1, 0, 0, 1, 0, 0, 0),
    @("Somalia", 1, 0, 0, 1, 0, 0, 0),
    @("Republica de Yibuti", 1, 0, 0, 1, 0, 0, 0),
    @("Papua Nueva Guinea", 1, 1, 0, 1, 0, 0, 0),
    @("Montserrat", 1, 0, 0, 1, 0, 0, 0),
    @("Fiyi", 1, 0, 0, 1, 0, 0, 0),
    @("Santa Sede", 1, 0, 0, 1, 0, 0, 0),
    @("Gambia", 1, 0, 0, 1, 0, 0, 0),
    @("San Vicente y las Granadinas", 1, 0, 0, 1, 0, 0, 0),
    @("Cabo Verde", 1, 1, 0, 1, 0, 0, 0),
    @("Nepal", 1, 0, 1, 0, 0, 0, 0),
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
}

$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 16:16"

